# Updates cryptos list figures (price/volume) per the GitHub Actions refresh.
# Price cells ("D" column) are prefixed with a literal leading apostrophe so Excel
# stores them as text (matching the source data, which uses "." both as a thousands
# separator and a decimal point) instead of silently reinterpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '''69.571.25'
$ws.Range("E2").Value = '  -0.25%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '''3.784.12'
$ws.Range("E3").Value = '  +0.75%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.02%  '

# Row 5: BNB
$ws.Range("D5").Value = '''614.03'
$ws.Range("E5").Value = '  -1.04%  '

# Row 6: Solana
$ws.Range("D6").Value = '''177.07'
$ws.Range("E6").Value = '  -1.71%  '

# Row 7: LidoStakedEther
$ws.Range("D7").Value = '''3.779.93'
$ws.Range("E7").Value = '  +0.71%  '

# Row 8: USDC
$ws.Range("E8").Value = '  -0.02%  '

# Row 9: XRP
$ws.Range("E9").Value = '  -0.66%  '

# Row 10: Dogecoin
$ws.Range("E10").Value = '  -1.65%  '

# Row 11: Toncoin
$ws.Range("D11").Value = '''6.41'
$ws.Range("E11").Value = '  +1.41%  '

# Row 12: Cardano
$ws.Range("E12").Value = '  -1.25%  '

# Row 13: Avalanche
$ws.Range("D13").Value = '''39.82'
$ws.Range("E13").Value = '  -3.35%  '

# Row 14: ShibaInu
$ws.Range("E14").Value = '  -2.13%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '''4.412.28'
$ws.Range("E15").Value = '  +0.83%  '

# Row 16: WrappedEther
$ws.Range("D16").Value = '''3.784.64'
$ws.Range("E16").Value = '  +0.59%  '

# Row 17: WrappedBTC
$ws.Range("D17").Value = '''69.595.62'
$ws.Range("E17").Value = '  -0.37%  '

# Row 18: Polkadot
$ws.Range("E18").Value = '  -1.14%  '

# Row 19: TRON
$ws.Range("E19").Value = '  -3.57%  '

# Row 20: ranking swapped - was BitcoinCash, now Chainlink
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '''16.59'
$ws.Range("E20").Value = '  -0.94%  '

# Row 21: ranking swapped - was Chainlink, now BitcoinCash
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '''508.84'
$ws.Range("E21").Value = '  +0.25%  '

# Row 22: Uniswap
$ws.Range("E22").Value = '  +0.50%  '

# Row 23: Polygon
$ws.Range("E23").Value = '  +0.70%  '

# Row 24: Fetch.AI
$ws.Range("E24").Value = '  -1.91%  '

# Row 25: Litecoin
$ws.Range("D25").Value = '''86.19'
$ws.Range("E25").Value = '  -1.23%  '

# Row 26: PEPE
$ws.Range("D26").Value = '''0.0000143'
$ws.Range("E26").Value = '  +3.92%  '

# Row 27: InternetComputer(DFINITY)
$ws.Range("D27").Value = '''12.85'
$ws.Range("E27").Value = '  -2.41%  '

# Row 28: RenderToken
$ws.Range("D28").Value = '''10.55'
$ws.Range("E28").Value = '  -5.47%  '

# Row 30: PancakeSwap
$ws.Range("D30").Value = '''2.99'
$ws.Range("E30").Value = '  +3.18%  '

# Row 31: ImmutableX
$ws.Range("E31").Value = '  -0.35%  '

# Row 32: NEARProtocol
$ws.Range("D32").Value = '''8.13'
$ws.Range("E32").Value = '  +2.82%  '

# Row 33: EthereumClassic
$ws.Range("D33").Value = '''31.24'
$ws.Range("E33").Value = '  +0.44%  '

# Row 34: Hedera
$ws.Range("E34").Value = '  -0.70%  '

# Row 35: FirstDigitalUSD
$ws.Range("E35").Value = '  -0.03%  '

# Row 36: Mantle
$ws.Range("D36").Value = '''1.04'
$ws.Range("E36").Value = '  -1.76%  '

# Row 37: Filecoin
$ws.Range("D37").Value = '''6.12'
$ws.Range("E37").Value = '  -1.50%  '

# Row 38: Kaspa
$ws.Range("E38").Value = '  +6.96%  '

# Row 39: Bittensor
$ws.Range("D39").Value = '''483.88'
$ws.Range("E39").Value = '  +13.12%  '

# Row 40: TheGraph
$ws.Range("D40").Value = '''0.340'
$ws.Range("E40").Value = '  +0.51%  '

# Row 41: Stacks
$ws.Range("E41").Value = '  -2.76%  '

# Row 42: OKB
$ws.Range("E42").Value = '  -0.83%  '

# Row 43: dogwifhat
$ws.Range("D43").Value = '''2.99'
$ws.Range("E43").Value = '  +4.14%  '

# Row 44: Arweave
$ws.Range("D44").Value = '''44.10'
$ws.Range("E44").Value = '  -3.79%  '

# Row 46: Maker
$ws.Range("D46").Value = '''2.940.31'
$ws.Range("E46").Value = '  -2.18%  '

# Row 47: VeChain
$ws.Range("E47").Value = '  -0.71%  '

# Row 48: InjectiveProtocol
$ws.Range("D48").Value = '''27.38'
$ws.Range("E48").Value = '  -0.16%  '

# Row 49: Monero
$ws.Range("D49").Value = '''139.27'
$ws.Range("E49").Value = '  +1.51%  '

# Row 51: ThetaToken
$ws.Range("E51").Value = '  -1.80%  '

Write-Host "Applied cryptos update"
